$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.017209748119836
$ws.Range("D2").Value = 1.019271901313267
$ws.Range("E2").Value = 1.01864176367692
$ws.Range("F2").Value = 1.015567788015579
$ws.Range("I2").Value = 1.027372539976271
$ws.Range("J2").Value = 1.022425163064438
$ws.Range("K2").Value = 1.022115641795462
$ws.Range("L2").Value = 1.021487374417678
$ws.Range("M2").Value = 1.01842255860826
$ws.Range("N2").Value = 1.011588584536039
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.018117239705064
$ws.Range("D3").Value = 1.020047674554548
$ws.Range("E3").Value = 1.019409418406767
$ws.Range("F3").Value = 1.017116733111366
$ws.Range("I3").Value = 1.027482867656249
$ws.Range("J3").Value = 1.022968662840647
$ws.Range("K3").Value = 1.022697487441704
$ws.Range("L3").Value = 1.022060988366195
$ws.Range("M3").Value = 1.019774634373723
$ws.Range("N3").Value = 1.011767981732139
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.018704030588423
$ws.Range("D4").Value = 1.02054958243443
$ws.Range("E4").Value = 1.019906196766855
$ws.Range("F4").Value = 1.018118504646188
$ws.Range("I4").Value = 1.027551869223494
$ws.Range("J4").Value = 1.02331933446014
$ws.Range("K4").Value = 1.023073258013221
$ws.Range("L4").Value = 1.022431554153625
$ws.Range("M4").Value = 1.020648547105376
$ws.Range("N4").Value = 1.011883714803986
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.018950617505497
$ws.Range("D5").Value = 1.02076056729115
$ws.Range("E5").Value = 1.020115054730839
$ws.Range("F5").Value = 1.018539535039388
$ws.Range("I5").Value = 1.027580305226257
$ws.Range("J5").Value = 1.023466514956072
$ws.Range("K5").Value = 1.023231058627128
$ws.Range("L5").Value = 1.022587195770308
$ws.Range("M5").Value = 1.021015712023315
$ws.Range("N5").Value = 1.011932285340558
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.018992014666277
$ws.Range("D6").Value = 1.020795991548824
$ws.Range("E6").Value = 1.020150123605488
$ws.Range("F6").Value = 1.018610221340328
$ws.Range("I6").Value = 1.027585046180945
$ws.Range("J6").Value = 1.023491213018257
$ws.Range("K6").Value = 1.023257543869015
$ws.Range("L6").Value = 1.022613320212136
$ws.Range("M6").Value = 1.021077347386189
$ws.Range("N6").Value = 1.011940435638784
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.018707325887942
$ws.Range("D7").Value = 1.020552401692883
$ws.Range("E7").Value = 1.019908987489593
$ws.Range("F7").Value = 1.018124130917802
$ws.Range("I7").Value = 1.027552251436122
$ws.Range("J7").Value = 1.023321302043829
$ws.Range("K7").Value = 1.023075367234335
$ws.Range("L7").Value = 1.022433634411835
$ws.Range("M7").Value = 1.020653454071421
$ws.Range("N7").Value = 1.011884364134706
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.017516525362105
$ws.Range("D8").Value = 1.019534091483676
$ws.Range("E8").Value = 1.018901184511215
$ws.Range("F8").Value = 1.016091367842981
$ws.Range("I8").Value = 1.027410319785108
$ws.Range("J8").Value = 1.02260905013259
$ws.Range("K8").Value = 1.022312428627296
$ws.Range("L8").Value = 1.021681354093789
$ws.Range("M8").Value = 1.018879702583796
$ws.Range("N8").Value = 1.011649284833646
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.015414996517432
$ws.Range("D9").Value = 1.017739194894601
$ws.Range("E9").Value = 1.017125752965868
$ws.Range("F9").Value = 1.0125053191128
$ws.Range("I9").Value = 1.027141956229966
$ws.Range("J9").Value = 1.021346254086508
$ws.Range("K9").Value = 1.020962515547299
$ws.Range("L9").Value = 1.020351155778459
$ws.Range("M9").Value = 1.015746493867288
$ws.Range("N9").Value = 1.01123237764821
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.014011839746396
$ws.Range("D10").Value = 1.016542291425214
$ws.Range("E10").Value = 1.015942468549628
$ws.Range("F10").Value = 1.010111559335011
$ws.Range("I10").Value = 1.026950802287338
$ws.Range("J10").Value = 1.02049921352103
$ws.Range("K10").Value = 1.020058881666227
$ws.Range("L10").Value = 1.019461295054372
$ws.Range("M10").Value = 1.013652277973579
$ws.Range("N10").Value = 1.010952651621902
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.013403750084476
$ws.Range("D11").Value = 1.0160239525268
$ws.Range("E11").Value = 1.0154301800454
$ws.Range("F11").Value = 1.009074234721516
$ws.Range("I11").Value = 1.026865132675841
$ws.Range("J11").Value = 1.020131209395441
$ws.Range("K11").Value = 1.019666725467469
$ws.Range("L11").Value = 1.019075252759025
$ws.Range("M11").Value = 1.012744113373002
$ws.Range("N11").Value = 1.010831103955377
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.013177801192734
$ws.Range("D12").Value = 1.015831407983697
$ws.Range("E12").Value = 1.015239905988242
$ws.Range("F12").Value = 1.00868879786645
$ws.Range("I12").Value = 1.026832876118061
$ws.Range("J12").Value = 1.019994331508602
$ws.Range("K12").Value = 1.019520929603845
$ws.Range("L12").Value = 1.018931750464087
$ws.Range("M12").Value = 1.012406571476346
$ws.Range("N12").Value = 1.010785892004997
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.013226271498551
$ws.Range("D13").Value = 1.015872709906215
$ws.Range("E13").Value = 1.01528071983727
$ws.Range("F13").Value = 1.008771481291009
$ws.Range("I13").Value = 1.026839814939094
$ws.Range("J13").Value = 1.020023700641007
$ws.Range("K13").Value = 1.019552209264279
$ws.Range("L13").Value = 1.018962537123048
$ws.Range("M13").Value = 1.012478984911316
$ws.Range("N13").Value = 1.010795593005488
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.013385074656836
$ws.Range("D14").Value = 1.016008036944709
$ws.Range("E14").Value = 1.015414451670711
$ws.Range("F14").Value = 1.00904237704593
$ws.Range("I14").Value = 1.026862475211061
$ws.Range("J14").Value = 1.020119898795791
$ws.Range("K14").Value = 1.019654676619593
$ws.Range("L14").Value = 1.019063393037066
$ws.Range("M14").Value = 1.012716216352673
$ws.Range("N14").Value = 1.010827368021455
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.013482908273474
$ws.Range("D15").Value = 1.016091415033118
$ws.Range("E15").Value = 1.0154968499639
$ws.Range("F15").Value = 1.009209267664557
$ws.Range("I15").Value = 1.026876379317317
$ws.Range("J15").Value = 1.020179145170512
$ws.Range("K15").Value = 1.019717792705742
$ws.Range("L15").Value = 1.019125519256361
$ws.Range("M15").Value = 1.012862354596397
$ws.Range("N15").Value = 1.010846937212902
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.014052185821504
$ws.Range("D16").Value = 1.016576690391736
$ws.Range("E16").Value = 1.01597646920673
$ws.Range("F16").Value = 1.010180385517048
$ws.Range("I16").Value = 1.026956426863581
$ws.Range("J16").Value = 1.020523610830841
$ws.Range("K16").Value = 1.020084889320903
$ws.Range("L16").Value = 1.019486900114879
$ws.Range("M16").Value = 1.013712520843103
$ws.Range("N16").Value = 1.010960709400533
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.014409140991636
$ws.Range("D17").Value = 1.016881071861645
$ws.Range("E17").Value = 1.016277343828188
$ws.Range("F17").Value = 1.010789320440194
$ws.Range("I17").Value = 1.027005862831586
$ws.Range("J17").Value = 1.020739355760002
$ws.Range("K17").Value = 1.020314924693301
$ws.Range("L17").Value = 1.019713390349109
$ws.Range("M17").Value = 1.014245440939837
$ws.Range("N17").Value = 1.011031962076
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.014617297367873
$ws.Range("D18").Value = 1.017058605470697
$ws.Range("E18").Value = 1.016452846733592
$ws.Range("F18").Value = 1.011144424082196
$ws.Range("I18").Value = 1.027034418284356
$ws.Range("J18").Value = 1.020865077452226
$ws.Range("K18").Value = 1.020449015830921
$ws.Range("L18").Value = 1.019845428140372
$ws.Range("M18").Value = 1.014556153820756
$ws.Range("N18").Value = 1.011073481600736
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.01468826491567
$ws.Range("D19").Value = 1.017119138652106
$ws.Range("E19").Value = 1.016512690002487
$ws.Range("F19").Value = 1.011265492248592
$ws.Range("I19").Value = 1.027044107498193
$ws.Range("J19").Value = 1.020907925170174
$ws.Range("K19").Value = 1.020494723089159
$ws.Range("L19").Value = 1.019890437713837
$ws.Range("M19").Value = 1.014662076933394
$ws.Range("N19").Value = 1.011087631737859
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.014370848204738
$ws.Range("D20").Value = 1.016848415326738
$ws.Range("E20").Value = 1.016245062015454
$ws.Range("F20").Value = 1.010723995578775
$ws.Range("I20").Value = 1.027000587748582
$ws.Range("J20").Value = 1.020716220645689
$ws.Range("K20").Value = 1.020290252824672
$ws.Range("L20").Value = 1.019689097342906
$ws.Range("M20").Value = 1.014188277183467
$ws.Range("N20").Value = 1.011024321574162
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.013338313220494
$ws.Range("D21").Value = 1.015968186788777
$ws.Range("E21").Value = 1.015375070624621
$ws.Range("F21").Value = 1.008962608592205
$ws.Range("I21").Value = 1.026855814331743
$ws.Range("J21").Value = 1.020091575937333
$ws.Range("K21").Value = 1.019624506168934
$ws.Range("L21").Value = 1.019033696491616
$ws.Range("M21").Value = 1.012646363423413
$ws.Range("N21").Value = 1.01081801283038
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.012688670074491
$ws.Range("D22").Value = 1.015414692196611
$ws.Range("E22").Value = 1.014828146221327
$ws.Range("F22").Value = 1.007854409099411
$ws.Range("I22").Value = 1.026762272686558
$ws.Range("J22").Value = 1.019697767998165
$ws.Range("K22").Value = 1.019205163535211
$ws.Range("L22").Value = 1.018620989416421
$ws.Range("M22").Value = 1.011675688339882
$ws.Range("N22").Value = 1.010687929596036
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.013033100667504
$ws.Range("D23").Value = 1.015708115711611
$ws.Range("E23").Value = 1.015118074057439
$ws.Range("F23").Value = 1.008441959469109
$ws.Range("I23").Value = 1.026812099305657
$ws.Range("J23").Value = 1.019906634365911
$ws.Range("K23").Value = 1.01942753712038
$ws.Range("L23").Value = 1.01883983295635
$ws.Range("M23").Value = 1.012190378367706
$ws.Range("N23").Value = 1.01075692412207
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.014388151219696
$ws.Range("D24").Value = 1.016863171430577
$ws.Range("E24").Value = 1.016259648752819
$ws.Range("F24").Value = 1.010753513310596
$ws.Range("I24").Value = 1.027002972195768
$ws.Range("J24").Value = 1.020726674774555
$ws.Range("K24").Value = 1.020301401241412
$ws.Range("L24").Value = 1.019700074523321
$ws.Range("M24").Value = 1.014214107425931
$ws.Range("N24").Value = 1.011027774114565
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.015958669134282
$ws.Range("D25").Value = 1.018203274931771
$ws.Range("E25").Value = 1.017584687873256
$ws.Range("F25").Value = 1.013432913770783
$ws.Range("I25").Value = 1.027213495276575
$ws.Range("J25").Value = 1.021673630356283
$ws.Range("K25").Value = 1.021312152406244
$ws.Range("L25").Value = 1.020695585034569
$ws.Range("M25").Value = 1.016557436442932
$ws.Range("N25").Value = 1.01134047366322
